# Add a new "2022-Q3" sheet right after "总计" (the summary sheet),
# populate it with fund-holdings data, and insert a matching new row
# into the summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after the first sheet.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Data rows — B..G are stored as text in the source file (so leading
# zeros in fund codes like "005662" survive), H as a real number.
$rows = @(
    @(0, "671010", "西部利得策略优选混合A", "3.33", "93.42", "8.94", "0.2977", 4),
    @(1, "005662", "嘉实金融精选股票A",     "4.94", "92.37", "5.12", "0.2529", 10),
    @(2, "005663", "嘉实金融精选股票C",     "2.60", "92.37", "5.12", "0.1331", 10),
    @(3, "011060", "西部利得策略优选混合C", "0.47", "93.42", "8.94", "0.0420", 4),
    @(4, "011124", "富国金融地产行业混合C", "0.90", "88.10", "4.00", "0.0360", 7),
    @(5, "006652", "富国金融地产行业混合A", "0.82", "88.10", "4.00", "0.0328", 7),
    @(6, "007674", "工银产业升级股票A",     "0.29", "94.12", "5.48", "0.0159", 8),
    @(7, "007675", "工银产业升级股票C",     "0.25", "94.12", "5.48", "0.0137", 8)
)

# Force columns B,D,E,F,G to be stored as text (not auto-coerced to
# numbers), matching the source data which keeps fund codes / percents
# / amounts as plain text.
$q3.Range("B2:B9").NumberFormat = "@"
$q3.Range("D2:G9").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $q3.Range("A$r").Value = $row[0]
    $q3.Range("B$r").Value = $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").Value = $row[3]
    $q3.Range("E$r").Value = $row[4]
    $q3.Range("F$r").Value = $row[5]
    $q3.Range("G$r").Value = $row[6]
    $q3.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# Match the look-and-feel of the other quarter sheets: bold, centered,
# bordered header row (B1:H1) and index column (A2:A9).
$headerAndIndex = @($q3.Range("B1:H1"), $q3.Range("A2:A9"))
foreach ($rng in $headerAndIndex) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
}

# ---------------------------------------------------------------
# 2. Insert a matching "2022-Q3" row into the summary sheet, right
#    after the header row, shifting the existing quarters down.
# ---------------------------------------------------------------
$summary.Rows("2:2").Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 8
$summary.Range("D2").Value = 0.82

# Rows("2:2").Insert() copies the header row's formatting down (bold
# for B:D, plain for A) — the opposite of the rest of the table (bold
# index column A, plain B:D). Put it back in line with the other rows.
$summary.Range("B2:D2").Font.Bold = $false
$summary.Range("B2:D2").Borders.LineStyle = -4142
$summary.Range("A2").Font.Bold = $true
$summary.Range("A2").HorizontalAlignment = -4108
$summary.Range("A2").VerticalAlignment = -4160
$summary.Range("A2").Borders.LineStyle = 1

# Renumber the index column (A) for the shifted rows so it keeps
# counting 0,1,2,... down the table.
for ($i = 3; $i -le 9; $i++) {
    $summary.Range("A$i").Value = $i - 2
}
